# preparing for scans with COMSOL data
# Add a new slide (the 4th) using the "Title and Content" layout (same
# layout used by the existing "Comparison of different gases" slide),
# with the title filled in and the content placeholder left blank so it
# is ready to receive the COMSOL gas-comparison scan data later.

$p = $ppt.ActivePresentation

# Layout index 2 on the slide master is "Title and Content".
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Comparison of different gases"
$titleRange.Font.Bold = $true
$titleRange.Font.Color.RGB = 6299648
